$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.322262525558472
$ws.Range("B1").Value = 2.377776622772217
$ws.Range("C1").Value = 2.950720310211182
$ws.Range("D1").Value = 3.396380662918091
$ws.Range("E1").Value = 1.670630693435669
